$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> [date (dash format), D, E, F, G, H]
$data = @{
    3  = @("28-07-2022", 1, 0, 0, 1, 1)
    4  = @("01-08-2022", 1, 1, 0, 0, 0)
    5  = @("04-08-2022", 1, 1, 0, 0, 0)
    6  = @("08-08-2022", 1, 1, 0, 0, 0)
    7  = @("11-08-2022", 0, 0, 0, 0, 1)
    8  = @("15-08-2022", 0, 0, 0, 0, 1)
    9  = @("18-08-2022", 0, 0, 0, 0, 1)
    10 = @("22-08-2022", 1, 1, 0, 0, 0)
    11 = @("25-08-2022", 1, 1, 0, 0, 0)
    12 = @("29-08-2022", 1, 1, 0, 0, 0)
    13 = @("01-09-2022", 1, 1, 0, 0, 0)
    14 = @("05-09-2022", 0, 0, 0, 0, 1)
    15 = @("08-09-2022", 1, 1, 0, 0, 0)
    16 = @("12-09-2022", 1, 1, 0, 0, 0)
    17 = @("15-09-2022", 0, 0, 0, 0, 1)
    18 = @("19-09-2022", 0, 0, 0, 0, 1)
    19 = @("22-09-2022", 0, 0, 0, 0, 1)
    20 = @("26-09-2022", 1, 1, 0, 0, 0)
    21 = @("29-09-2022", 0, 0, 0, 0, 1)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]

    # Force column A to remain plain text so Excel doesn't reinterpret the
    # dash-separated date string as a date value, then restore the default
    # (unstyled) cell formatting so the style matches the original workbook.
    $aCell = $ws.Range("A$r")
    $aCell.NumberFormat = "@"
    $aCell.Value = $vals[0]
    $aCell.Style = "Normal"

    $ws.Range("D$r").Value = $vals[1]
    $ws.Range("E$r").Value = $vals[2]
    $ws.Range("F$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
    $ws.Range("H$r").Value = $vals[5]
}
